$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the login sample values: username field changes, password stays the same text
$ws.Range("A2").Value = "visulonsprint"
$ws.Range("B2").Value = "admin@123"

# Move the active selection to C16 (blank cell) as captured in the sheet view
$ws.Range("C16").Select()
